# Updates "Price" (D) and "Volume(1h)" (E) columns for the cryptos table.
# All D/E cells hold text (not numbers), so plain numeric-looking strings
# (e.g. "20.80", "0.3550") are entered with a leading apostrophe where
# needed to force a text entry and keep trailing zeros verbatim, matching
# the source data exactly instead of being auto-coerced to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.543.78'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '1.754.33'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '324.27'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = '0.4555'
$ws.Range("E7").Value = '  +1.73%  '
$ws.Range("D8").Value = "'0.3550"
$ws.Range("E8").Value = '  -1.92%  '
$ws.Range("D9").Value = '0.07471'
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("D10").Value = '41.51'
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("D11").Value = '1.087'
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = "'20.80"
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("D14").Value = '6.016'
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("D15").Value = '7.171'
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("D16").Value = '1.754.57'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = "'93.90"
$ws.Range("E17").Value = '  +1.13%  '
$ws.Range("D18").Value = '0.00001055'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").Value = '0.06395'
$ws.Range("E19").Value = '  -0.47%  '
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '17.11'
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").Value = '5.743'
$ws.Range("E22").Value = '  -1.66%  '
$ws.Range("D23").Value = '27.598.77'
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").Value = '2.084'
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("E26").Value = '  +2.07%  '
$ws.Range("D27").Value = '20.11'
$ws.Range("E27").Value = '  -1.45%  '
$ws.Range("D28").Value = '1.961.44'
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("D29").Value = '2.136'
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("D30").Value = '125.65'
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("D31").Value = '1.091'
$ws.Range("E31").Value = '  +0.78%  '
$ws.Range("D32").Value = "'0.09220"
$ws.Range("E32").Value = '  +2.15%  '
$ws.Range("D33").Value = '3.658'
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("D34").Value = '5.533'
$ws.Range("E34").Value = '  -0.29%  '
$ws.Range("D35").Value = '11.72'
$ws.Range("E35").Value = '  -2.42%  '
$ws.Range("D36").Value = '0.02282'
$ws.Range("E36").Value = '  -1.76%  '
$ws.Range("D37").Value = '0.2096'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").Value = '0.06017'
$ws.Range("D39").Value = "'0.6290"
$ws.Range("E39").Value = '  -1.23%  '
$ws.Range("D40").Value = "'4.920"
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("E41").Value = '  -3.00%  '
$ws.Range("D42").Value = '1.389'
$ws.Range("E42").Value = '  +0.19%  '
$ws.Range("D43").Value = '7.808'
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").Value = '13.19'
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").Value = '3.717'
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").Value = '0.5861'
$ws.Range("E46").Value = '  -0.44%  '
$ws.Range("D47").Value = '122.14'
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("D49").Value = '0.06892'
$ws.Range("E49").Value = '  +0.30%  '
$ws.Range("D50").Value = '1.131'
$ws.Range("E50").Value = '  -2.64%  '
$ws.Range("D51").Value = '72.14'
$ws.Range("E51").Value = '  -0.39%  '
